# Fix the "5 a 14" age-group header typo ("5a 14" -> "5 a 14") on the
# per-virus detail sheets (Ad, Parainfluenza, Inf A, Inf B, Metapnemovirus).
# The "Global semana epidemiologica" and "VRS" sheets keep the original text.

$wb = $excel.ActiveWorkbook

$sheetNames = @("Ad", "Parainfluenza", "Inf A", "Inf B", "Metapnemovirus")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F1").Value = "5 a 14"
    $ws.Range("N1").Value = "5 a 14"
}
